# shifting time frame for FG1 heating season
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D14").Value = 10.7667291222218
$ws.Range("E14").Value = 2.980704309444402
$ws.Range("F14").Value = 0.2444383333333209
$ws.Range("G14").Value = 3.225142642777744
